$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would be auto-parsed as a number by Excel ---
# Pre-format as Text ("@") so the literal string is preserved exactly,
# matching the source workbook, where these columns store formatted
# price strings (e.g. "595.00") as text, not numeric values.
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "595.00"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "167.28"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.363"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "27.64"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "8.03"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "356.91"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "4.68"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "1.94"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "1.00"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "10.31"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "69.93"
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "547.61"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "158.14"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "2.43"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "152.78"

# --- Cells whose new text is safely non-numeric already ---
$ws.Cells.Item(2,4).Value = "67.462.00"
$ws.Cells.Item(2,5).Value = "  -0.40%  "
$ws.Cells.Item(3,4).Value = "2.634.42"
$ws.Cells.Item(4,5).Value = "  -0.01%  "
$ws.Cells.Item(5,5).Value = "  -0.80%  "
$ws.Cells.Item(6,5).Value = "  -0.19%  "
$ws.Cells.Item(8,5).Value = "  -2.31%  "
$ws.Cells.Item(9,4).Value = "2.633.85"
$ws.Cells.Item(9,5).Value = "  -1.51%  "
$ws.Cells.Item(10,5).Value = "  -2.82%  "
$ws.Cells.Item(11,5).Value = "  +1.32%  "
$ws.Cells.Item(12,5).Value = "  +0.02%  "
$ws.Cells.Item(14,5).Value = "  -1.13%  "
$ws.Cells.Item(15,4).Value = "3.113.45"
$ws.Cells.Item(15,5).Value = "  -1.42%  "
$ws.Cells.Item(16,5).Value = "  -1.85%  "
$ws.Cells.Item(17,4).Value = "67.384.87"
$ws.Cells.Item(17,5).Value = "  -0.43%  "
$ws.Cells.Item(18,4).Value = "2.627.77"
$ws.Cells.Item(18,5).Value = "  -1.48%  "
$ws.Cells.Item(19,5).Value = "  +1.78%  "
$ws.Cells.Item(20,5).Value = "  +2.87%  "
$ws.Cells.Item(21,5).Value = "  -2.02%  "
$ws.Cells.Item(22,5).Value = "  -1.82%  "
$ws.Cells.Item(23,5).Value = "  -3.13%  "
$ws.Cells.Item(24,2).Value = "SuiNetwork"
$ws.Cells.Item(24,3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(24,5).Value = "  -5.13%  "
$ws.Cells.Item(25,2).Value = "Dai"
$ws.Cells.Item(25,3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(25,5).Value = "  +0.05%  "
$ws.Cells.Item(26,5).Value = "  +1.15%  "
$ws.Cells.Item(27,5).Value = "  -1.34%  "
$ws.Cells.Item(28,5).Value = "  -1.70%  "
$ws.Cells.Item(29,5).Value = "  +0.16%  "
$ws.Cells.Item(30,5).Value = "  -1.83%  "
$ws.Cells.Item(31,5).Value = "  -1.64%  "
$ws.Cells.Item(32,5).Value = "  -1.23%  "
$ws.Cells.Item(33,5).Value = "  -3.17%  "
$ws.Cells.Item(34,5).Value = "  -2.04%  "
$ws.Cells.Item(35,5).Value = "  +4.46%  "
$ws.Cells.Item(37,5).Value = "  -3.92%  "
$ws.Cells.Item(38,5).Value = "  +1.62%  "
$ws.Cells.Item(39,5).Value = "  -2.79%  "
$ws.Cells.Item(40,5).Value = "  -1.80%  "
$ws.Cells.Item(41,5).Value = "  +2.01%  "
$ws.Cells.Item(42,5).Value = "  -1.17%  "
$ws.Cells.Item(43,5).Value = "  -1.73%  "
$ws.Cells.Item(45,5).Value = "  -3.84%  "
$ws.Cells.Item(46,4).Value = "0.0₆0300"
$ws.Cells.Item(46,5).Value = "  -0.51%  "
$ws.Cells.Item(47,5).Value = "  -0.61%  "
$ws.Cells.Item(48,5).Value = "  -1.95%  "
$ws.Cells.Item(50,5).Value = "  -1.63%  "
$ws.Cells.Item(51,5).Value = "  -0.74%  "
